$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B96").Value = "rcolker@nibs.org is`na valid deliverable e-mail box address."
$ws.Range("B97").Value = "pkmakin@Brilliant-Books.net is`na valid deliverable e-mail box address."
$ws.Range("B98").Value = "sara.grochowski@brilliant-books.net is`na valid deliverable e-mail box address."
$ws.Range("B99").Value = "It was not possible to determine if mitch.simpson@brokawsupply.com`nis a valid deliverable e-mail box address."
$ws.Range("B100").Value = "mbrown@brownrichards.com is not`na valid deliverable e-mail box address."

$ws.Rows.Item(96).EntireRow.AutoFit()
$ws.Rows.Item(97).EntireRow.AutoFit()
$ws.Rows.Item(98).EntireRow.AutoFit()
$ws.Rows.Item(99).EntireRow.AutoFit()
$ws.Rows.Item(100).EntireRow.AutoFit()
